$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Left table (A:E) ---
# Row 2 headers already present: A2 in: floor type, B2 in: length, C2 in: width, D2 out: room cost

# Row 3 - hardwood
$ws.Range("A3").Value = "hardwood"
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 11
$ws.Range("D3").Value = 152.9

# Row 4 - carpet
$ws.Range("A4").Value = "carpet"
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 12
$ws.Range("D4").Value = 335.16

# Row 5 - carpet
$ws.Range("A5").Value = "carpet"
$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 383.04

# Row 6 - tile
$ws.Range("A6").Value = "tile"
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 299.39999999999998

# Row 7 - tile
$ws.Range("A7").Value = "tile"
$ws.Range("B7").Value = 9
$ws.Range("C7").Value = 9
$ws.Range("D7").Value = 404.19

# Row 8 - final overall cost
$ws.Range("D8").Value = 1504.69

# --- Right lookup table (Table2, G1:H4) ---
$ws.Range("G2").Value = "hardwood"
$ws.Range("H2").Value = 1.39

$ws.Range("G3").Value = "carpet"
$ws.Range("H3").Value = 3.99

$ws.Range("G4").Value = "tile"
$ws.Range("H4").Value = 4.99

# --- Selection change ---
$ws.Range("D9").Select() | Out-Null
